# Generate Report for Archive
#
# The localization status moved on from "Ready for handoff" to
# "In Translation" everywhere it is reported (the Overview sheet's
# per-language status columns, and the "Status" column of each
# language detail sheet). The status columns are then narrowed to
# match the shorter text, mirroring how the reporting tool re-sizes
# columns whenever it regenerates this sheet.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: the zh-cn / de-de status columns (E & F) ---
foreach ($addr in "E2", "F2", "E3", "F3") {
    $cell = $wsOverview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- zh-cn / de-de detail sheets: the "Status" column (C) ---
foreach ($addr in "C2", "C3") {
    $cellZh = $wsZhCn.Range($addr)
    if ($cellZh.Value2 -eq $oldStatus) {
        $cellZh.Value = $newStatus
    }

    $cellDe = $wsDeDe.Range($addr)
    if ($cellDe.Value2 -eq $oldStatus) {
        $cellDe.Value = $newStatus
    }
}

# --- Narrow the status columns to fit the shorter "In Translation" text ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
